$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") values get updated from 45203 to 45204 for all data rows (2-261)
$ws.Range("C2:C261").Value = 45204
